# Update the crypto price/volume table with the latest scraped values.
# For D-column cells whose new text looks like a plain number (e.g. "1.001",
# "0.5250"), the cell is pre-formatted as Text ("@") so Excel stores the
# literal string (preserving trailing zeros / exact digits) instead of
# silently coercing it to a numeric value.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '26.489.74'
$ws.Range("E2").Value = '  +0.33%  '

# Row 3
$ws.Range("D3").Value = '1.836.47'
$ws.Range("E3").Value = '  -0.11%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.001'
$ws.Range("E4").Value = '  -0.02%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '259.64'
$ws.Range("E5").Value = '  -0.30%  '

# Row 6
$ws.Range("E6").Value = '  -0.02%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5250'

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3190'
$ws.Range("E8").Value = '  -2.05%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06772'
$ws.Range("E9").Value = '  +0.18%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '18.65'
$ws.Range("E10").Value = '  -0.26%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.7805'
$ws.Range("E11").Value = '  +2.09%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07736'
$ws.Range("E12").Value = '  +0.89%  '

# Row 13
$ws.Range("D13").Value = '1.832.24'
$ws.Range("E13").Value = '  -0.46%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '87.54'
$ws.Range("E14").Value = '  -1.53%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.006'
$ws.Range("E15").Value = '  -0.47%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.001'
$ws.Range("E16").Value = '  -0.14%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '13.78'
$ws.Range("E17").Value = '  -1.43%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '1.000'
$ws.Range("E18").Value = '  -0.01%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.000007914'
$ws.Range("E19").Value = '  +0.50%  '

# Row 20
$ws.Range("D20").Value = '26.505.68'
$ws.Range("E20").Value = '  +0.29%  '

# Row 21
$ws.Range("D21").Value = '2.068.57'
$ws.Range("E21").Value = '  -0.23%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.596'
$ws.Range("E22").Value = '  +0.74%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.362'
$ws.Range("E23").Value = '  -0.92%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '5.951'
$ws.Range("E24").Value = '  +0.16%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '141.39'
$ws.Range("E25").Value = '  -2.00%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.164'
$ws.Range("E26").Value = '  -3.24%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.673'
$ws.Range("E27").Value = '  +2.01%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '16.92'
$ws.Range("E28").Value = '  -0.33%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '111.51'
$ws.Range("E29").Value = '  +0.02%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.140'
$ws.Range("E30").Value = '  -0.82%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.08677'
$ws.Range("E31").Value = '  -0.28%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.047'
$ws.Range("E32").Value = '  -2.29%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.04849'
$ws.Range("E33").Value = '  +0.96%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.7257'
$ws.Range("E34").Value = '  +4.29%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.126'
$ws.Range("E35").Value = '  +0.17%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.851'
$ws.Range("E36").Value = '  -0.09%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.089'
$ws.Range("E37").Value = '  +0.64%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.228'
$ws.Range("E38").Value = '  +1.46%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01761'
$ws.Range("E39").Value = '  -0.11%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.4742'
$ws.Range("E40").Value = '  -1.88%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.8900'
$ws.Range("E41").Value = '  -0.21%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '109.36'
$ws.Range("E42").Value = '  -1.17%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.898'
$ws.Range("E43").Value = '  -3.16%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.001'
$ws.Range("E44").Value = '  -0.01%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '7.609'
$ws.Range("E45").Value = '  -0.99%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.4119'
$ws.Range("E46").Value = '  -0.94%  '

# Row 47
$ws.Range("B47").Value = 'Cronos'
$ws.Range("C47").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.05854'
$ws.Range("E47").Value = '  -0.14%  '

# Row 48
$ws.Range("B48").Value = 'EnergySwap'
$ws.Range("C48").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '8.942'
$ws.Range("E48").Value = '  -0.13%  '

# Row 49
$ws.Range("B49").Value = 'Elrond'
$ws.Range("C49").Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '34.76'
$ws.Range("E49").Value = '  -0.10%  '

# Row 50
$ws.Range("B50").Value = 'Algorand'
$ws.Range("C50").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.1221'
$ws.Range("E50").Value = '  -1.60%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.8905'
$ws.Range("E51").Value = '  +0.94%  '
